$d = $word.ActiveDocument

# The first paragraph in the body holds the ID placeholder text.
$p = $d.Paragraphs.Item(1)

# Replace the old ID token (plus the trailing space run) with the new ID token,
# collapsing the paragraph down to a single run.
$d.Content.Find.Execute("**ID__AFFARS_pgi_5327_topic_4__ID** ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_AFMC_PGI_5327_9001__ID**", 2)

# Update the paragraph formatting: new left indent and a (near-invisible) paragraph border.
$p.Format.LeftIndent = 11.25
$p.Format.Borders.DistanceFromTop = 5
$p.Format.Borders.DistanceFromLeft = 5
$p.Format.Borders.DistanceFromBottom = 5
$p.Format.Borders.DistanceFromRight = 5
